$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 40000
$ws.Cells.Item(95, 10).Value = 40000
$ws.Cells.Item(95, 12).Value = 40000
$ws.Cells.Item(95, 14).Value = -45492
$ws.Cells.Item(98, 8).Value = 980.9286
$ws.Cells.Item(98, 9).Value = 375.18182
$ws.Cells.Item(98, 11).Value = 375.18182
$ws.Cells.Item(98, 13).Value = 1122.81818
$ws.Cells.Item(122, 8).Value = 980.9286
$ws.Cells.Item(122, 9).Value = 375.18182
$ws.Cells.Item(122, 11).Value = 1125.54546
$ws.Cells.Item(122, 13).Value = 1324.45454
$ws.Cells.Item(138, 8).Value = 2966.125
$ws.Cells.Item(138, 9).Value = 1864.5
$ws.Cells.Item(138, 10).Value = 3333.3333
$ws.Cells.Item(138, 11).Value = 5593.5
$ws.Cells.Item(138, 12).Value = 9999.999899999999
$ws.Cells.Item(138, 13).Value = -453.5
$ws.Cells.Item(138, 14).Value = -20279.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2618.75
$ws.Cells.Item(63, 9).Value = 2158.3333
$ws.Cells.Item(63, 11).Value = 2158.3333
$ws.Cells.Item(63, 13).Value = -1472.3333
$ws.Cells.Item(66, 8).Value = 2618.75
$ws.Cells.Item(66, 9).Value = 2158.3333
$ws.Cells.Item(66, 11).Value = 10791.6665
$ws.Cells.Item(66, 13).Value = -7359.666499999999
$ws.Cells.Item(88, 8).Value = 2897.3845
$ws.Cells.Item(88, 9).Value = 993.4
$ws.Cells.Item(88, 11).Value = 993.4
$ws.Cells.Item(88, 13).Value = -587.4
$ws.Cells.Item(91, 8).Value = 2897.3845
$ws.Cells.Item(91, 9).Value = 993.4
$ws.Cells.Item(91, 11).Value = 993.4
$ws.Cells.Item(91, 13).Value = 410.6
$ws.Cells.Item(98, 8).Value = 40000
$ws.Cells.Item(98, 10).Value = 40000
$ws.Cells.Item(98, 12).Value = 40000
$ws.Cells.Item(98, 14).Value = -45990
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 34023.363
$ws.Cells.Item(82, 9).Value = 18419
$ws.Cells.Item(82, 10).Value = 39875
$ws.Cells.Item(82, 11).Value = 18419
$ws.Cells.Item(82, 12).Value = 39875
$ws.Cells.Item(82, 13).Value = -18036
$ws.Cells.Item(82, 14).Value = -40641
$ws.Cells.Item(85, 8).Value = 34023.363
$ws.Cells.Item(85, 9).Value = 18419
$ws.Cells.Item(85, 10).Value = 39875
$ws.Cells.Item(85, 11).Value = 18419
$ws.Cells.Item(85, 12).Value = 39875
$ws.Cells.Item(85, 13).Value = -17093
$ws.Cells.Item(85, 14).Value = -42527
$ws.Cells.Item(86, 8).Value = 2211.25
$ws.Cells.Item(86, 9).Value = 2230
$ws.Cells.Item(86, 10).Value = 2200
$ws.Cells.Item(86, 11).Value = 2230
$ws.Cells.Item(86, 12).Value = 2200
$ws.Cells.Item(86, 13).Value = -1107
$ws.Cells.Item(86, 14).Value = -4446
$ws.Cells.Item(89, 8).Value = 2211.25
$ws.Cells.Item(89, 9).Value = 2230
$ws.Cells.Item(89, 10).Value = 2200
$ws.Cells.Item(89, 11).Value = 11150
$ws.Cells.Item(89, 12).Value = 11000
$ws.Cells.Item(89, 13).Value = -5534
$ws.Cells.Item(89, 14).Value = -22232
$ws.Cells.Item(95, 8).Value = 15717
$ws.Cells.Item(95, 10).Value = 15717
$ws.Cells.Item(95, 12).Value = 15717
$ws.Cells.Item(95, 14).Value = -21209
$ws.Cells.Item(105, 8).Value = 2200
$ws.Cells.Item(105, 9).Value = 1900
$ws.Cells.Item(105, 11).Value = 1900
$ws.Cells.Item(105, 13).Value = -153
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 60142
$ws.Cells.Item(9, 10).Value = 60142
$ws.Cells.Item(9, 12).Value = 60142
$ws.Cells.Item(9, 14).Value = -60478
$ws.Cells.Item(16, 8).Value = 1000
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 13).Value = -713
$ws.Cells.Item(59, 8).Value = 34400
$ws.Cells.Item(59, 9).Value = 32000
$ws.Cells.Item(59, 11).Value = 32000
$ws.Cells.Item(59, 13).Value = -30855
$ws.Cells.Item(92, 8).Value = 32417.818
$ws.Cells.Item(92, 10).Value = 32709.6
$ws.Cells.Item(92, 12).Value = 32709.6
$ws.Cells.Item(92, 14).Value = -37701.6
$ws.Cells.Item(113, 8).Value = 1000
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(132, 8).Value = 1148.1538
$ws.Cells.Item(132, 9).Value = 1148.1538
$ws.Cells.Item(132, 11).Value = 3444.4614
$ws.Cells.Item(132, 13).Value = -914.4614000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 750
$ws.Cells.Item(5, 9).Value = 750
$ws.Cells.Item(5, 11).Value = 2250
$ws.Cells.Item(5, 13).Value = -2138
$ws.Cells.Item(36, 8).Value = 344.6
$ws.Cells.Item(36, 10).Value = 500
$ws.Cells.Item(36, 12).Value = 1500
$ws.Cells.Item(36, 14).Value = -1838
$ws.Cells.Item(135, 8).Value = 750
$ws.Cells.Item(135, 9).Value = 750
$ws.Cells.Item(135, 11).Value = 6750
$ws.Cells.Item(135, 13).Value = -4215
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 25250
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 25250
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 25250
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(57, 14).Value = -26890
$ws.Cells.Item(80, 8).Value = 5098.8
$ws.Cells.Item(80, 10).Value = 5624.75
$ws.Cells.Item(80, 12).Value = 5624.75
$ws.Cells.Item(80, 14).Value = -7620.75
$ws.Cells.Item(83, 8).Value = 5098.8
$ws.Cells.Item(83, 10).Value = 5624.75
$ws.Cells.Item(83, 12).Value = 28123.75
$ws.Cells.Item(83, 14).Value = -38107.75
$ws.Cells.Item(92, 8).Value = 8189.2
$ws.Cells.Item(92, 10).Value = 8189.2
$ws.Cells.Item(92, 12).Value = 8189.2
$ws.Cells.Item(92, 14).Value = -11933.2
$ws.Cells.Item(132, 8).Value = 3006
$ws.Cells.Item(132, 9).Value = 12
$ws.Cells.Item(132, 10).Value = 6000
$ws.Cells.Item(132, 11).Value = 36
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = 2494
$ws.Cells.Item(132, 14).Value = -23060
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1600
$ws.Cells.Item(93, 10).Value = 1600
$ws.Cells.Item(93, 12).Value = 1600
$ws.Cells.Item(93, 14).Value = -4096
$ws.Cells.Item(94, 8).Value = 49703.223
$ws.Cells.Item(94, 10).Value = 49703.223
$ws.Cells.Item(94, 12).Value = 49703.223
$ws.Cells.Item(94, 14).Value = -51055.223
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 9666.666999999999
$ws.Cells.Item(94, 8).Value = 5999.5
$ws.Cells.Item(94, 10).Value = 5999.5
$ws.Cells.Item(94, 12).Value = 5999.5
$ws.Cells.Item(94, 14).Value = -7801.5
$ws.Cells.Item(97, 8).Value = 29357.666
$ws.Cells.Item(97, 10).Value = 29357.666
$ws.Cells.Item(97, 12).Value = 29357.666
$ws.Cells.Item(97, 14).Value = -31339.666
$ws.Cells.Item(132, 8).Value = 2020.6
$ws.Cells.Item(132, 9).Value = 2020.6
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 6061.799999999999
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -3531.799999999999
$ws.Cells.Item(132, 14).ClearContents()
